# Update the arithmetic problems in the worksheet table.
# Each data row in the table (rows 1, 5, 9, 13, 17) has 5 cells with a
# "NN÷N=" style expression. Replace each cell's text in place (leaving
# the trailing paragraph/cell mark untouched) so run formatting such as
# rFonts/sz is preserved.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = $newText
}

# Row 1
Set-CellText $t 1 1 "31÷4="
Set-CellText $t 1 2 "53÷8="
Set-CellText $t 1 3 "18÷9="
Set-CellText $t 1 4 "20÷3="
Set-CellText $t 1 5 "49÷2="

# Row 5
Set-CellText $t 5 1 "49÷3="
Set-CellText $t 5 2 "65÷8="
Set-CellText $t 5 3 "50÷2="
Set-CellText $t 5 4 "11÷9="
Set-CellText $t 5 5 "67÷5="

# Row 9
Set-CellText $t 9 1 "80÷4="
Set-CellText $t 9 2 "29÷8="
Set-CellText $t 9 3 "95÷9="
Set-CellText $t 9 4 "54÷3="
Set-CellText $t 9 5 "93÷5="

# Row 13
Set-CellText $t 13 1 "95÷2="
Set-CellText $t 13 2 "70÷7="
Set-CellText $t 13 3 "63÷2="
Set-CellText $t 13 4 "27÷2="
Set-CellText $t 13 5 "99÷2="

# Row 17
Set-CellText $t 17 1 "35÷5="
Set-CellText $t 17 2 "19÷3="
Set-CellText $t 17 3 "51÷4="
Set-CellText $t 17 4 "30÷6="
Set-CellText $t 17 5 "99÷9="
